$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row before row 4 for the new task description
#    (formatting is inherited from the row above, matching row 3's style)
$ws.Rows("4:4").Insert()

# 2) Fill the new row 4 with the new task text and style it like rows 2/3
$ws.Range("A4:C4").Merge()
$ws.Range("A4:C4").Value = "Beheben der Probleme von unterschiedlichen Dateiversionen"
$ws.Rows("4:4").RowHeight = 69
$ws.Range("A4:C4").HorizontalAlignment = -4108
$ws.Range("A4:C4").VerticalAlignment = -4160
$ws.Range("A4:C4").WrapText = $true

$wb.Save()
